$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.424.85'
$ws.Range('E2').Value = '  -2.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.889.14'
$ws.Range('E3').Value = '  -1.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.70'
$ws.Range('E5').Value = '  -4.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.73'
$ws.Range('E6').Value = '  -2.86%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.888.67'
$ws.Range('E9').Value = '  -1.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.04'
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.147'
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.03'
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.366.40'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.435.77'
$ws.Range('E17').Value = '  -2.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.55'
$ws.Range('E18').Value = '  -2.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.884.71'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '432.03'
$ws.Range('E20').Value = '  -2.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  -2.23%  '
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.84'
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.96'
$ws.Range('E25').Value = '  +1.44%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('E27').Value = '  -11.15%  '
$ws.Range('E28').Value = '  -5.96%  '
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('E30').Value = '  -2.63%  '
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('E32').Value = '  -7.12%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -1.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.46'
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.961'
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.39'
$ws.Range('E37').Value = '  -3.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.89'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.93'
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('E40').Value = '  -10.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.21'
$ws.Range('E41').Value = '  -3.42%  '
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.14'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('E44').Value = '  -4.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.712.01'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '132.89'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0335'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '343.61'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.58'
$ws.Range('E51').Value = '  -5.10%  '
